$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Forgotten UserName" (row 7) as fixed
$ws.Range("C7").Value = "x"
$ws.Range("D7").Value = Get-Date -Year 2017 -Month 6 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Mark "ResetPassword" (row 13) as fixed
$ws.Range("C13").Value = "x"
$ws.Range("D13").Value = Get-Date -Year 2017 -Month 6 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
